$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3403.5833
$ws.Range("I40").Value = 2995.5
$ws.Range("K40").Value = 2995.5
$ws.Range("M40").Value = -2820.5
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H111").Value = 766.63635
$ws.Range("I111").Value = 766.63635
$ws.Range("K111").Value = 2299.90905
$ws.Range("M111").Value = 767.0909499999998
$ws.Range("H113").Value = 5092.154
$ws.Range("I113").Value = 2856.8572
$ws.Range("J113").Value = 7700
$ws.Range("K113").Value = 2856.8572
$ws.Range("L113").Value = 7700
$ws.Range("M113").Value = 397.1428000000001
$ws.Range("N113").Value = -14208
$ws.Range("H132").Value = 200008000
$ws.Range("I132").Value = 250008500
$ws.Range("K132").Value = 750025500
$ws.Range("M132").Value = -750022970
$ws.Range("H137").Value = 2387.25
$ws.Range("I137").Value = 1780.2858
$ws.Range("J137").Value = 3237
$ws.Range("K137").Value = 5340.857400000001
$ws.Range("L137").Value = 9711
$ws.Range("M137").Value = -2790.857400000001
$ws.Range("N137").Value = -14811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5343.5537
$ws.Range("I32").Value = 4531.618
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 4531.618
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -4244.618
$ws.Range("N32").Value = -50574
$ws.Range("H86").Value = 200000
$ws.Range("J86").Value = 200000
$ws.Range("L86").Value = 200000
$ws.Range("N86").Value = -202372
$ws.Range("H89").Value = 200000
$ws.Range("J89").Value = 200000
$ws.Range("L89").Value = 600000
$ws.Range("N89").Value = -611856
$ws.Range("H104").Value = 12612.25
$ws.Range("J104").Value = 12612.25
$ws.Range("L104").Value = 12612.25
$ws.Range("N104").Value = -19600.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 37089.57
$ws.Range("J100").Value = 37089.57
$ws.Range("L100").Value = 37089.57
$ws.Range("N100").Value = -39253.57
$ws.Range("H105").Value = 2989.9
$ws.Range("I105").Value = 2179.8
$ws.Range("J105").Value = 3800
$ws.Range("K105").Value = 2179.8
$ws.Range("L105").Value = 3800
$ws.Range("M105").Value = -432.8000000000002
$ws.Range("N105").Value = -7294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4190.5
$ws.Range("I31").Value = 3187.75
$ws.Range("J31").Value = 5527.5
$ws.Range("K31").Value = 3187.75
$ws.Range("L31").Value = 5527.5
$ws.Range("M31").Value = -2892.75
$ws.Range("N31").Value = -6117.5
$ws.Range("H34").Value = 4190.5
$ws.Range("I34").Value = 3187.75
$ws.Range("J34").Value = 5527.5
$ws.Range("K34").Value = 3187.75
$ws.Range("L34").Value = 5527.5
$ws.Range("M34").Value = -2985.75
$ws.Range("N34").Value = -5931.5
$ws.Range("H43").Value = 6569
$ws.Range("J43").Value = 6569
$ws.Range("L43").Value = 6569
$ws.Range("N43").Value = -6937
$ws.Range("H86").Value = 4722.8184
$ws.Range("I86").Value = 4851.857
$ws.Range("J86").Value = 4497
$ws.Range("K86").Value = 4851.857
$ws.Range("L86").Value = 4497
$ws.Range("M86").Value = -3728.857
$ws.Range("N86").Value = -6743
$ws.Range("H89").Value = 4722.8184
$ws.Range("I89").Value = 4851.857
$ws.Range("J89").Value = 4497
$ws.Range("K89").Value = 24259.285
$ws.Range("L89").Value = 22485
$ws.Range("M89").Value = -18643.285
$ws.Range("N89").Value = -33717
$ws.Range("H101").Value = 6569
$ws.Range("J101").Value = 6569
$ws.Range("L101").Value = 6569
$ws.Range("N101").Value = -13059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1329.8334
$ws.Range("I80").Value = 1396.8
$ws.Range("K80").Value = 4190.4
$ws.Range("M80").Value = -3254.4
$ws.Range("H83").Value = 1329.8334
$ws.Range("I83").Value = 1396.8
$ws.Range("K83").Value = 12571.2
$ws.Range("M83").Value = -7891.199999999999
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H133").Value = 1500
$ws.Range("I133").Value = 1500
$ws.Range("K133").Value = 4500
$ws.Range("M133").Value = 560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 216.85715
$ws.Range("I107").Value = 228
$ws.Range("K107").Value = 228
$ws.Range("M107").Value = 1692
$ws.Range("H113").Value = 2057
$ws.Range("I113").Value = 1776.25
$ws.Range("J113").Value = 2244.1667
$ws.Range("K113").Value = 1776.25
$ws.Range("L113").Value = 2244.1667
$ws.Range("M113").Value = 393.75
$ws.Range("N113").Value = -6584.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3607.7273
$ws.Range("I7").Value = 3685.75
$ws.Range("K7").Value = 3685.75
$ws.Range("M7").Value = -3573.75
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H99").Value = 41749.5
$ws.Range("I99").Value = 41749.5
$ws.Range("K99").Value = 41749.5
$ws.Range("M99").Value = -38754.5
$ws.Range("H126").Value = 3607.7273
$ws.Range("I126").Value = 3685.75
$ws.Range("K126").Value = 11057.25
$ws.Range("M126").Value = -8587.25
$ws.Range("H136").Value = 3621
$ws.Range("I136").Value = 3443.2
$ws.Range("J136").Value = 3917.3333
$ws.Range("K136").Value = 10329.6
$ws.Range("L136").Value = 11751.9999
$ws.Range("M136").Value = -7779.599999999999
$ws.Range("N136").Value = -16851.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 10400.667
$ws.Range("J101").Value = 10400.667
$ws.Range("L101").Value = 10400.667
$ws.Range("N101").Value = -16890.667
